# Apply updated cryptocurrency price/volume figures (GitHub Actions data refresh)

function Set-TextValue($range, $value) {
    # Force the assigned value to be stored as literal text (preserving the
    # original cell style) so numeric-looking strings such as "1.003" are not
    # auto-converted into numbers by Excel.
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-39: update Price (D) and Volume(1h) (E) values
Set-TextValue $ws.Range("D2") "24.628.43"
Set-TextValue $ws.Range("E2") "  -0.24%  "
Set-TextValue $ws.Range("D3") "1.689.47"
Set-TextValue $ws.Range("E3") "  +0.05%  "
Set-TextValue $ws.Range("D4") "1.003"
Set-TextValue $ws.Range("E4") "  +0.06%  "
Set-TextValue $ws.Range("D5") "313.36"
Set-TextValue $ws.Range("E5") "  -0.68%  "
Set-TextValue $ws.Range("D6") "1.003"
Set-TextValue $ws.Range("E6") "  +0.07%  "
Set-TextValue $ws.Range("D7") "0.3896"
Set-TextValue $ws.Range("E7") "  -1.08%  "
Set-TextValue $ws.Range("D8") "0.4028"
Set-TextValue $ws.Range("E8") "  -0.67%  "
Set-TextValue $ws.Range("D9") "1.498"
Set-TextValue $ws.Range("E9") "  +0.80%  "
Set-TextValue $ws.Range("D10") "1.002"
Set-TextValue $ws.Range("E10") "  -0.15%  "
Set-TextValue $ws.Range("D11") "52.80"
Set-TextValue $ws.Range("E11") "  +0.16%  "
Set-TextValue $ws.Range("D12") "0.08729"
Set-TextValue $ws.Range("E12") "  -1.21%  "
Set-TextValue $ws.Range("D13") "7.613"
Set-TextValue $ws.Range("E13") "  +5.11%  "
Set-TextValue $ws.Range("D14") "24.84"
Set-TextValue $ws.Range("E14") "  +5.59%  "
Set-TextValue $ws.Range("D15") "7.957"
Set-TextValue $ws.Range("E15") "  -0.72%  "
Set-TextValue $ws.Range("D16") "0.00001347"
Set-TextValue $ws.Range("E16") "  +2.46%  "
Set-TextValue $ws.Range("D17") "1.683.92"
Set-TextValue $ws.Range("E17") "  -0.30%  "
Set-TextValue $ws.Range("D18") "98.18"
Set-TextValue $ws.Range("E18") "  -1.42%  "
Set-TextValue $ws.Range("D19") "0.07092"
Set-TextValue $ws.Range("E19") "  +1.18%  "
Set-TextValue $ws.Range("D20") "19.74"
Set-TextValue $ws.Range("E20") "  +1.21%  "
Set-TextValue $ws.Range("D21") "7.278"
Set-TextValue $ws.Range("E21") "  +4.05%  "
Set-TextValue $ws.Range("D22") "1.003"
Set-TextValue $ws.Range("E22") "  -0.34%  "
Set-TextValue $ws.Range("D23") "14.21"
Set-TextValue $ws.Range("E23") "  -0.57%  "
Set-TextValue $ws.Range("D24") "24.622.28"
Set-TextValue $ws.Range("E24") "  -0.19%  "
Set-TextValue $ws.Range("D25") "3.000"
Set-TextValue $ws.Range("E25") "  -8.71%  "
Set-TextValue $ws.Range("D26") "2.346"
Set-TextValue $ws.Range("E26") "  -0.45%  "
Set-TextValue $ws.Range("D27") "22.69"
Set-TextValue $ws.Range("E27") "  -0.10%  "
Set-TextValue $ws.Range("D28") "161.43"
Set-TextValue $ws.Range("E28") "  -0.77%  "
Set-TextValue $ws.Range("D29") "8.574"
Set-TextValue $ws.Range("E29") "  +12.52%  "
Set-TextValue $ws.Range("D30") "5.220"
Set-TextValue $ws.Range("E30") "  +0.63%  "
Set-TextValue $ws.Range("D31") "136.24"
Set-TextValue $ws.Range("E31") "  +0.70%  "
Set-TextValue $ws.Range("D32") "1.869.02"
Set-TextValue $ws.Range("E32") "  -0.53%  "
Set-TextValue $ws.Range("D33") "0.08760"
Set-TextValue $ws.Range("E33") "  +2.63%  "
Set-TextValue $ws.Range("D34") "7.402"
Set-TextValue $ws.Range("E34") "  +4.15%  "
Set-TextValue $ws.Range("D35") "1.037"
Set-TextValue $ws.Range("E35") "  -2.13%  "
Set-TextValue $ws.Range("D36") "1.986"
Set-TextValue $ws.Range("E36") "  +5.27%  "
Set-TextValue $ws.Range("D37") "0.02910"
Set-TextValue $ws.Range("E37") "  +7.34%  "
Set-TextValue $ws.Range("D38") "0.2710"
Set-TextValue $ws.Range("E38") "  -0.62%  "
Set-TextValue $ws.Range("D39") "10.77"
Set-TextValue $ws.Range("E39") "  -4.41%  "

# Rows 40-41: Stellar and Aptos swapped positions, with updated Price/Volume values
Set-TextValue $ws.Range("B40") "Aptos"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D40") "14.18"
Set-TextValue $ws.Range("E40") "  -1.73%  "

Set-TextValue $ws.Range("B41") "Stellar"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D41") "0.09104"
Set-TextValue $ws.Range("E41") "  -0.80%  "

# Rows 42-51: update Price (D) and Volume(1h) (E) values
Set-TextValue $ws.Range("D42") "0.7788"
Set-TextValue $ws.Range("E42") "  +2.15%  "
Set-TextValue $ws.Range("D43") "1.454"
Set-TextValue $ws.Range("E43") "  -0.58%  "
Set-TextValue $ws.Range("D44") "16.60"
Set-TextValue $ws.Range("E44") "  +4.26%  "
Set-TextValue $ws.Range("D45") "0.7176"
Set-TextValue $ws.Range("E45") "  +0.55%  "
Set-TextValue $ws.Range("D46") "2.576"
Set-TextValue $ws.Range("E46") "  -0.41%  "
Set-TextValue $ws.Range("D47") "4.196"
Set-TextValue $ws.Range("E47") "  -0.56%  "
Set-TextValue $ws.Range("D48") "1.003"
Set-TextValue $ws.Range("E48") "  +0.07%  "
Set-TextValue $ws.Range("D49") "1.335"
Set-TextValue $ws.Range("E49") "  +1.27%  "
Set-TextValue $ws.Range("D50") "137.62"
Set-TextValue $ws.Range("E50") "  -1.42%  "
Set-TextValue $ws.Range("D51") "90.67"
Set-TextValue $ws.Range("E51") "  +1.20%  "
